# Insert a new price record as row 179 in the "Feria Lagunitas de Puerto Montt -
# Zapallo italiano" sheet, pushing the existing rows 179-248 down to 180-249.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(179).Insert()

$ws.Range("A179").Value = 4
$ws.Range("B179").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C179").Value = "Los Lagos"
$ws.Range("D179").Value = 44726
$ws.Range("E179").Value = 10
$ws.Range("F179").Value = 100112032
$ws.Range("G179").Value = "Zapallo italiano"
$ws.Range("H179").Value = "Sin especificar"
$ws.Range("I179").Value = "Primera"
$ws.Range("J179").Value = 200
$ws.Range("K179").Value = 15000
$ws.Range("L179").Value = 15000
$ws.Range("M179").Value = 15000
$ws.Range("N179").Value = "$/caja 50 unidades"
$ws.Range("O179").Value = "Región de Arica y Parinacota"
$ws.Range("P179").Value = 300
$ws.Range("Q179").Value = 50
$ws.Range("R179").Value = "Hortaliza"
